# Update "想去人数" (want-to-go count, column F) figures to the freshly
# scraped values, across the "展览", "演出" and "全部类型" sheets
# (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

$exhibition = $wb.Worksheets.Item("展览")
$exhibition.Range("F2").Value  = 1201
$exhibition.Range("F3").Value  = 2014
$exhibition.Range("F5").Value  = 1290
$exhibition.Range("F9").Value  = 362
$exhibition.Range("F10").Value = 145
$exhibition.Range("F12").Value = 897
$exhibition.Range("F19").Value = 722
$exhibition.Range("F24").Value = 938
$exhibition.Range("F25").Value = 391
$exhibition.Range("F26").Value = 211
$exhibition.Range("F28").Value = 321
$exhibition.Range("F31").Value = 439

$performance = $wb.Worksheets.Item("演出")
$performance.Range("F7").Value  = 270
$performance.Range("F11").Value = 137
$performance.Range("F12").Value = 30

$allTypes = $wb.Worksheets.Item("全部类型")
$allTypes.Range("F2").Value  = 338
$allTypes.Range("F3").Value  = 1201
$allTypes.Range("F4").Value  = 2014
$allTypes.Range("F6").Value  = 1290
$allTypes.Range("F11").Value = 362
$allTypes.Range("F12").Value = 145
$allTypes.Range("F14").Value = 897
$allTypes.Range("F24").Value = 270
$allTypes.Range("F26").Value = 722
$allTypes.Range("F31").Value = 938
$allTypes.Range("F32").Value = 391
$allTypes.Range("F35").Value = 211
$allTypes.Range("F37").Value = 321
$allTypes.Range("F39").Value = 137
$allTypes.Range("F42").Value = 30
$allTypes.Range("F43").Value = 439

$localLife = $wb.Worksheets.Item("本地生活")
$localLife.Range("F2").Value = 338
